$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 317; this shifts the former rows 317-366 down to 318-367
$ws.Rows(317).Insert()

# Populate the newly inserted row 317 with the new record's data
$ws.Range("A317").Value2 = 2
$ws.Range("B317").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C317").Value2 = "Coquimbo"
$ws.Range("D317").Value2 = 44951
$ws.Range("E317").Value2 = 4
$ws.Range("F317").Value2 = 100112021
$ws.Range("G317").Value2 = "Ají"
$ws.Range("H317").Value2 = "Americana (o)"
$ws.Range("I317").Value2 = "Primera"
$ws.Range("J317").Value2 = 240
$ws.Range("K317").Value2 = 13000
$ws.Range("L317").Value2 = 14000
$ws.Range("M317").Value2 = 13500
$ws.Range("N317").Value2 = "$/caja 25 kilos"
$ws.Range("O317").Value2 = "Provincia de Limarí"
$ws.Range("P317").Value2 = 540
$ws.Range("Q317").Value2 = 25
$ws.Range("R317").Value2 = "Hortaliza"
